# Add a new "Sheet2" after "Sheet1" with a small 3x3 sample grid of names,
# a conditional format highlighting "Mark", a (now-orphaned) hidden
# _FilterDatabase defined name left over from an earlier AutoFilter/table
# experiment, and update the view/selection state on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New worksheet placed immediately after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Row 3 - fill right-to-left (D, C, B) to match the original shared-string
# insertion order (Tim, Craig, Mark, ...).
$ws2.Range("D3").Value = "Tim"
$ws2.Range("C3").Value = "Craig"
$ws2.Range("B3").Value = "Mark"

# Row 5 - filled next.
$ws2.Range("B5").Value = "Mark"
$ws2.Range("C5").Value = "Tim"
$ws2.Range("D5").Value = "Pat"

# Row 4 - filled last (and given a custom row height), matching the
# later-inserted-in-the-middle feel of the original edit.
$ws2.Range("B4").Value = "Vinnie"
$ws2.Range("C4").Value = "Boom "
$ws2.Range("D4").Value = "Bop"
$ws2.Rows.Item(4).RowHeight = 12.75

# Conditional formatting: highlight cells equal to "Mark" (Excel's
# built-in "Light Red Fill with Dark Red Text" look).
$cf = $ws2.Range("B3:D5").FormatConditions.Add(3, 3, "=""Mark""")
$cf.Font.Color = 393372
$cf.Interior.Color = 13551615

# Leftover hidden _FilterDatabase name scoped to Sheet2 (vestige of an
# earlier AutoFilter/table on a larger A5:F49 range).
$fd = $ws2.Names.Add("_xlnm._FilterDatabase", "=Sheet2!`$A`$5:`$F`$49")
$fd.Visible = $false

# Print orientation.
$ws2.PageSetup.Orientation = 1

# Selections / active sheet.
$ws1.Range("E10").Select() | Out-Null
$ws2.Range("G6").Select() | Out-Null
